$d = $word.ActiveDocument

# 20 simple one-to-one cell text replacements (unique strings, safe via Find/Replace)
$d.Content.Find.Execute("37÷2=18, 1", $true, $true, $false, $false, $false, $true, 1, $false, "12÷5=2, 2", 2) | Out-Null
$d.Content.Find.Execute("86÷7=12, 2", $true, $true, $false, $false, $false, $true, 1, $false, "20÷4=5, 0", 2) | Out-Null
$d.Content.Find.Execute("66÷3=22, 0", $true, $true, $false, $false, $false, $true, 1, $false, "74÷4=18, 2", 2) | Out-Null
$d.Content.Find.Execute("65÷4=16, 1", $true, $true, $false, $false, $false, $true, 1, $false, "57÷8=7, 1", 2) | Out-Null
$d.Content.Find.Execute("83÷5=16, 3", $true, $true, $false, $false, $false, $true, 1, $false, "78÷3=26, 0", 2) | Out-Null
$d.Content.Find.Execute("48÷9=5, 3", $true, $true, $false, $false, $false, $true, 1, $false, "30÷5=6, 0", 2) | Out-Null
$d.Content.Find.Execute("66÷4=16, 2", $true, $true, $false, $false, $false, $true, 1, $false, "21÷8=2, 5", 2) | Out-Null
$d.Content.Find.Execute("64÷8=8, 0", $true, $true, $false, $false, $false, $true, 1, $false, "69÷9=7, 6", 2) | Out-Null
$d.Content.Find.Execute("22÷5=4, 2", $true, $true, $false, $false, $false, $true, 1, $false, "44÷9=4, 8", 2) | Out-Null
$d.Content.Find.Execute("54÷2=27, 0", $true, $true, $false, $false, $false, $true, 1, $false, "49÷6=8, 1", 2) | Out-Null
$d.Content.Find.Execute("90÷9=10, 0", $true, $true, $false, $false, $false, $true, 1, $false, "49÷3=16, 1", 2) | Out-Null
$d.Content.Find.Execute("75÷8=9, 3", $true, $true, $false, $false, $false, $true, 1, $false, "75÷3=25, 0", 2) | Out-Null
$d.Content.Find.Execute("15÷2=7, 1", $true, $true, $false, $false, $false, $true, 1, $false, "85÷2=42, 1", 2) | Out-Null
$d.Content.Find.Execute("12÷8=1, 4", $true, $true, $false, $false, $false, $true, 1, $false, "47÷6=7, 5", 2) | Out-Null
$d.Content.Find.Execute("66÷9=7, 3", $true, $true, $false, $false, $false, $true, 1, $false, "23÷9=2, 5", 2) | Out-Null
$d.Content.Find.Execute("42÷8=5, 2", $true, $true, $false, $false, $false, $true, 1, $false, "29÷8=3, 5", 2) | Out-Null
$d.Content.Find.Execute("76÷2=38, 0", $true, $true, $false, $false, $false, $true, 1, $false, "18÷7=2, 4", 2) | Out-Null
$d.Content.Find.Execute("45÷8=5, 5", $true, $true, $false, $false, $false, $true, 1, $false, "83÷7=11, 6", 2) | Out-Null
$d.Content.Find.Execute("29÷3=9, 2", $true, $true, $false, $false, $false, $true, 1, $false, "22÷9=2, 4", 2) | Out-Null
$d.Content.Find.Execute("81÷4=20, 1", $true, $true, $false, $false, $false, $true, 1, $false, "90÷7=12, 6", 2) | Out-Null

# Row 9 (the "76÷3=25, 1 / 29÷4=7, 1 / 20÷2=10, 0 / 38÷7=5, 3 / 21÷5=4, 1" row) is re-laid out:
# a new cell is inserted after column 2 and the last cell is dropped, net effect keeps 5
# cells with identical formatting -> apply as direct per-cell text assignment.
$t = $d.Tables.Item(1)
$t.Cell(9, 1).Range.Text = "35÷6=5, 5"
$t.Cell(9, 2).Range.Text = "82÷9=9, 1"
$t.Cell(9, 3).Range.Text = "87÷9=9, 6"
$t.Cell(9, 4).Range.Text = "20÷2=10, 0"
$t.Cell(9, 5).Range.Text = "95÷3=31, 2"

Write-Output "edit complete"
